# Ajout d'un titre pour le site
#
# 1) Fix a typo on the existing 'meta description vide' row (B6): drop the
#    curly quotes around the word "description".
# 2) Add a new audit row (row 7) for the missing <title> meta tag, in the
#    same shape as the existing rows (Categorie/Probleme/Explication/Bonne
#    pratique/Action recommandee/Reference incl. hyperlink).
# 3) Move the active-cell selection from D8 to A8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) fix the typo in B6 ---
$ws.Range("B6").Value = "(index.html l.6): meta description vide"

# --- 2) new row 7 ---
$ws.Range("A7").Value = "SEO/accessibilité"
$ws.Range("B7").Value = "(index.html l.22): meta title vide"
$ws.Range("C7").Value = "la meta titre est ce qui sera affiché en premier lors de l’apparition du site web dans les résultats de recherche et se doit donc d’être attractif"
$ws.Range("D7").Value = "utiliser un titre concis et reprenant des mots clés de manière pertinente et naturelle"
$ws.Range("E7").Value = "exemple de titre : “La Chouette agence | Experts web design à Lyon”"
$ws.Range("F7").Value = "https://university.webflow.com/lesson/seo-title-meta-description?utm_source=google&utm_medium=search&utm_campaign=general-paid-workhorse&utm_term=keyword-targeting&utm_content=dynamic-search-ads-webflow-university-t1&gclid=Cj0KCQjw7MGJBhD-ARIsAMZ0eevWyJeq"

# Hyperlink on F7 pointing at the webflow article (same url as the display text)
$url = "https://university.webflow.com/lesson/seo-title-meta-description?utm_source=google&utm_medium=search&utm_campaign=general-paid-workhorse&utm_term=keyword-targeting&utm_content=dynamic-search-ads-webflow-university-t1&gclid=Cj0KCQjw7MGJBhD-ARIsAMZ0eevWyJeq"
$null = $wb.Hyperlinks.Add($ws.Range("F7"), $url, [System.Type]::Missing, [System.Type]::Missing, $url)

# Adding the hyperlink auto-applies Excel's built-in blue/underlined
# 'Hyperlink' style; the other reference cell (F6) keeps the plain row
# formatting, so reset F7 back to the same plain format used across the row.
$ws.Range("A7").Copy()
$null = $ws.Range("F7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# --- 3) move the active selection to A8 (was D8) ---
$null = $ws.Range("A8").Select()
